$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.223.59'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '3.332.52'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'577.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("D6").Value = "'183.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.56%  '

$ws.Range("D9").Value = "'0.128"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("D12").Value = '3.914.36'
$ws.Range("E12").Value = '  +2.30%  '

$ws.Range("D14").Value = "'27.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").Value = '67.427.42'
$ws.Range("E15").Value = '  -0.85%  '

$ws.Range("E16").Value = '  -0.19%  '

$ws.Range("D17").Value = '3.333.24'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").Value = "'444.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.21%  '

$ws.Range("E19").Value = '  +2.59%  '

$ws.Range("D20").Value = "'5.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").Value = "'7.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.80%  '

$ws.Range("D22").Value = "'73.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.26%  '

$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").Value = '3.488.95'

$ws.Range("E25").Value = '  +1.02%  '

$ws.Range("D26").Value = "'0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.00%  '

$ws.Range("E27").Value = '  +3.29%  '

$ws.Range("E28").Value = '  -2.62%  '

$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("E30").Value = '  +1.65%  '

$ws.Range("D31").Value = "'22.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.41%  '

$ws.Range("D32").Value = "'5.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.27%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = "'6.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.78%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = "'1.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("D38").Value = "'27.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.87%  '

$ws.Range("E39").Value = '  -1.64%  '

$ws.Range("D40").Value = '2.830.97'
$ws.Range("E40").Value = '  +8.12%  '

$ws.Range("D41").Value = "'0.790"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").Value = "'4.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.86%  '

$ws.Range("D43").Value = "'6.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("E44").Value = '  -0.54%  '

$ws.Range("D45").Value = "'0.0669"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.22%  '

$ws.Range("D46").Value = "'24.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '

$ws.Range("E47").Value = '  -2.43%  '

$ws.Range("D48").Value = "'321.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.50%  '

$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("D50").Value = "'0.981"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").Value = "'30.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.60%  '
